$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.154.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.294.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.88%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.90%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.286.16"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.94%  "

$ws.Range("E10").Value = "  -3.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.89%  "

$ws.Range("E13").Value = "  -4.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.811.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "17.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.69%  "

$ws.Range("E17").Value = "  -4.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.289.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.050.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.962"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "415.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.60%  "

$ws.Range("E28").Value = "  -5.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.80%  "

$ws.Range("E31").Value = "  -2.65%  "

$ws.Range("E32").Value = "  -3.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "571.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.07%  "

$ws.Range("E34").Value = "  -4.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("E37").Value = "  -1.54%  "

$ws.Range("E38").Value = "  +4.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0734"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.53%  "

$ws.Range("E41").Value = "  -5.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.091.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.26%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.99%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.07%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.98%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0397"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.44%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.127"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.78%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.73%  "
